$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.36
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.73
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.62
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 19
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 17
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 7.5
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("AA5").Value = 23
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 13
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 1.91
$ws.Range("Z5").Value = 23
$ws.Range("M6").Value = 1.1
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63
$ws.Range("AZ7").Value = 41
$ws.Range("G7").Value = 2.63
$ws.Range("I7").Value = 2.55
$ws.Range("L7").Value = 3.2
$ws.Range("Q7").Value = 1.93
$ws.Range("R7").Value = 1.93
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("AC9").Value = 12
$ws.Range("AE9").Value = 13
$ws.Range("AI9").Value = 10
$ws.Range("AK9").Value = 17
$ws.Range("G9").Value = 3.6
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 1.85
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 2.25
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 12
$ws.Range("O9").Value = 1.19
$ws.Range("AB10").Value = 27
$ws.Range("AC10").Value = 14.5
$ws.Range("AD10").Value = 11.25
$ws.Range("AE10").Value = 25
$ws.Range("AH10").Value = 40
$ws.Range("AI10").Value = 150
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 2.95
$ws.Range("AO10").Value = 4.6
$ws.Range("AP10").Value = 15.5
$ws.Range("AQ10").Value = 10.25
$ws.Range("AR10").Value = 35
$ws.Range("AS10").Value = 200
$ws.Range("AT10").Value = 3.45
$ws.Range("AU10").Value = 10
$ws.Range("AV10").Value = 90
$ws.Range("AW10").Value = 15.5
$ws.Range("AX10").Value = 110
$ws.Range("G10").Value = 1.14
$ws.Range("H10").Value = 5.7
$ws.Range("I10").Value = 17
$ws.Range("J10").Value = 1.5
$ws.Range("K10").Value = 2.67
$ws.Range("L10").Value = 11.75
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 5.24
$ws.Range("T10").Value = 3.91
$ws.Range("U10").Value = 2.33
$ws.Range("V10").Value = 1.56
$ws.Range("W10").Value = 6.5
$ws.Range("X10").Value = 5
$ws.Range("Y10").Value = 8.5
$ws.Range("Z10").Value = 5.4
$ws.Range("N11").Value = 17
$ws.Range("O11").Value = 1.17
$ws.Range("P11").Value = 5
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 2.35
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 21
$ws.Range("AB14").Value = 26
$ws.Range("AC14").Value = 26
$ws.Range("AT14").Value = 4.33
$ws.Range("G14").Value = 4
$ws.Range("I14").Value = 1.7
$ws.Range("S14").Value = 1.2
$ws.Range("T14").Value = 4.33
$ws.Range("X14").Value = 29
$ws.Range("Q15").Value = 1.65
$ws.Range("R15").Value = 2.2
$ws.Range("AE16").Value = 17
$ws.Range("AK16").Value = 51
$ws.Range("AT16").Value = 3.25
$ws.Range("Q16").Value = 1.75
$ws.Range("R16").Value = 2.05
$ws.Range("S16").Value = 1.33
$ws.Range("T16").Value = 3.25
$ws.Range("M17").Value = 1.02
$ws.Range("O17").Value = 1.13
$ws.Range("M18").Value = 1.05
$ws.Range("O18").Value = 1.29
$ws.Range("M19").Value = 1.04
$ws.Range("O19").Value = 1.2
$ws.Range("M20").Value = 1.02
$ws.Range("O20").Value = 1.14
$ws.Range("O21").Value = 1.08
$ws.Range("U22").Value = 1.57
$ws.Range("V23").Value = 1.73
$ws.Range("M24").Value = 1.05
$ws.Range("N24").Value = 8
$ws.Range("O24").Value = 1.37
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("M25").Value = 1.03
$ws.Range("O25").Value = 1.22
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("J26").Value = 2.88
$ws.Range("M26").Value = 1.03
$ws.Range("O26").Value = 1.19
$ws.Range("U26").Value = 1.62
$ws.Range("AA27").Value = 17
$ws.Range("AG27").Value = 351
$ws.Range("AH27").Value = 10
$ws.Range("AN27").Value = 3.75
$ws.Range("AO27").Value = 11
$ws.Range("AP27").Value = 23
$ws.Range("AX27").Value = 23
$ws.Range("G27").Value = 1.91
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 4.2
$ws.Range("J27").Value = 2.75
$ws.Range("K27").Value = 2.05
$ws.Range("L27").Value = 4.5
$ws.Range("N27").Value = 8
$ws.Range("O27").Value = 1.33
$ws.Range("P27").Value = 3
$ws.Range("Q27").Value = 2.25
$ws.Range("R27").Value = 1.62
$ws.Range("X27").Value = 8.5
$ws.Range("Y27").Value = 9
$ws.Range("Z27").Value = 17
$ws.Range("N32").Value = 15
$ws.Range("Q32").Value = 1.6
$ws.Range("R32").Value = 2.3
$ws.Range("M35").Value = 1.07
$ws.Range("O35").Value = 1.36
$ws.Range("M36").Value = 1.05
$ws.Range("O36").Value = 1.29
$ws.Range("M37").Value = 1.06
$ws.Range("N37").Value = 8
$ws.Range("M38").Value = 1.08
$ws.Range("O38").Value = 1.44
$ws.Range("P38").Value = 2.63
$ws.Range("AE39").Value = 13
$ws.Range("AF39").Value = 41
$ws.Range("AJ39").Value = 11
$ws.Range("AL39").Value = 23
$ws.Range("AU39").Value = 7.5
$ws.Range("AW39").Value = 5
$ws.Range("AX39").Value = 17
$ws.Range("BA39").Value = 67
$ws.Range("G39").Value = 2.2
$ws.Range("H39").Value = 3.4
$ws.Range("I39").Value = 3.1
$ws.Range("X39").Value = 11
$ws.Range("Z39").Value = 21
